# "Images of results and tables" — strip the worked examples (polynomial
# regression / matrix definition rows) from the "Writing" sheet, leaving
# just the empty (but still formatted) rows behind, drop the now-dead
# hyperlinks, widen column A, and move the active selection down to A7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two Wikipedia hyperlinks (regresión lineal / mínimos cuadrados)
# that lived on A10:A11 — their backing text is being removed too.
$ws.Range("A10").Hyperlinks.Delete()
$ws.Range("A11").Hyperlinks.Delete()

# Clear the text of rows 9-14 (the polynomial-regression / pizza / matrix
# example content) while keeping each row's existing cell formatting.
$ws.Range("A9:A14").ClearContents()

# Column A is noticeably wider in the new layout.
$ws.Columns.Item(1).ColumnWidth = 58

# Active cell/selection moves from E2 to A7.
$ws.Range("A7").Select()
